$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A23").NumberFormat = "@"
$ws.Range("A23").Value = "2026-02-26"
$ws.Range("B23").Value = "Parcialmente Nublado"
$ws.Range("C23").Value = "Abafado"
$ws.Range("D23").Value = "normal"
$ws.Range("E23").Value = "aula"
$ws.Range("F23").Value = "nenhuma"
$ws.Range("G23").Value = 20
$ws.Range("H23").Value = 31.5
$ws.Range("I23").Value = 2
$ws.Range("J23").Value = "jtq01"
$ws.Range("K23").Value = "Terça e Quinta"
$ws.Range("L23").Value = "08:00"
$ws.Range("M23").Value = "Jefferson"
